$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.461.31'
$ws.Range("E2").Value = '  +0.98%  '

$ws.Range("D3").Value = '2.237.56'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("E4").Value = '  +1.35%  '

$ws.Range("D5").Value = '''306.86'
$ws.Range("E5").Value = '  +0.56%  '

$ws.Range("D6").Value = '''93.77'
$ws.Range("E6").Value = '  -2.12%  '

$ws.Range("D7").Value = '''0.571'
$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("E8").Value = '  +0.20%  '

$ws.Range("E9").Value = '  -0.64%  '

$ws.Range("D10").Value = '''34.60'
$ws.Range("E10").Value = '  -1.15%  '

$ws.Range("D11").Value = '''0.0801'
$ws.Range("E11").Value = '  -0.83%  '

$ws.Range("E12").Value = '  -0.56%  '

$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("D14").Value = '''0.831'
$ws.Range("E14").Value = '  +0.43%  '

$ws.Range("D15").Value = '2.201.20'
$ws.Range("E15").Value = '  -1.70%  '

$ws.Range("D16").Value = '''13.50'
$ws.Range("E16").Value = '  -0.44%  '

$ws.Range("D17").Value = '44.080.14'
$ws.Range("E17").Value = '  +0.40%  '

$ws.Range("E18").Value = '  -0.78%  '

$ws.Range("D19").Value = '''6.34'
$ws.Range("E19").Value = '  +1.74%  '

$ws.Range("D20").Value = '''11.89'
$ws.Range("E20").Value = '  -3.33%  '

$ws.Range("D21").Value = '''65.68'
$ws.Range("E21").Value = '  +1.40%  '

$ws.Range("D22").Value = '''237.74'
$ws.Range("E22").Value = '  +0.59%  '

$ws.Range("E23").Value = '  +1.20%  '

$ws.Range("E24").Value = '  +0.93%  '

$ws.Range("E25").Value = '  -0.28%  '

$ws.Range("E26").Value = '  +3.90%  '

$ws.Range("D27").Value = '''9.76'
$ws.Range("E27").Value = '  -1.82%  '

$ws.Range("E28").Value = '  -0.59%  '

$ws.Range("D29").Value = '''5.90'
$ws.Range("E29").Value = '  -0.98%  '

$ws.Range("D30").Value = '''19.92'
$ws.Range("E30").Value = '  -0.57%  '

$ws.Range("D31").Value = '''153.65'
$ws.Range("E31").Value = '  -1.16%  '

$ws.Range("E32").Value = '  -1.86%  '

$ws.Range("E33").Value = '  +0.43%  '

$ws.Range("D34").Value = '''3.10'
$ws.Range("E34").Value = '  -5.55%  '

$ws.Range("E35").Value = '  +2.28%  '

$ws.Range("E36").Value = '  +0.63%  '

$ws.Range("D37").Value = '''1.78'
$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").Value = '''14.83'
$ws.Range("E38").Value = '  -3.49%  '

$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("D40").Value = '''3.76'
$ws.Range("E40").Value = '  -1.68%  '

$ws.Range("E41").Value = '  -1.03%  '

$ws.Range("E42").Value = '  +0.29%  '

$ws.Range("D43").Value = '1.769.63'
$ws.Range("E43").Value = '  +1.81%  '

$ws.Range("D44").Value = '''0.192'
$ws.Range("E44").Value = '  +1.66%  '

$ws.Range("D45").Value = '''78.89'
$ws.Range("E45").Value = '  -7.67%  '

$ws.Range("D46").Value = '''98.56'
$ws.Range("E46").Value = '  -1.51%  '

$ws.Range("D47").Value = '''4.87'
$ws.Range("E47").Value = '  -1.30%  '

$ws.Range("D48").Value = '''69.69'
$ws.Range("E48").Value = '  +0.52%  '

$ws.Range("D49").Value = '''8.08'
$ws.Range("E49").Value = '  -0.14%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '''1.57'
$ws.Range("E50").Value = '  +3.82%  '

$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '''54.57'
$ws.Range("E51").Value = '  +0.33%  '
